# Add data for 2022-10-15 (carjacking-by-neighborhood-by-month)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet and update header label to reflect the new "through" date.
$ws.Name = "Through 2022-10-07"
$ws.Range("B1").Value = "October 2022 (through October 07)"

# Row 2 - Garfield Park
$ws.Range("L2").Value = 4
$ws.Range("V2").Value = 6

# Row 8 - West Loop
$ws.Range("B8").Value = 1

# Row 10 - Grand Crossing
$ws.Range("L10").Value = 3

# Row 13 - Roseland
$ws.Range("BJ13").Value = 1

# Row 17 - Washington Heights
$ws.Range("L17").Value = 2

# Row 23 - Auburn Gresham
$ws.Range("AP23").Value = 1

# Row 24 - Riverdale
$ws.Range("L24").Value = 1

# Row 45 - Bridgeport
$ws.Range("B45").Value = 1

# Row 47 - Bucktown
$ws.Range("V47").Value = 1

# Row 48 - Calumet Heights
$ws.Range("L48").Value = 2

# Row 66 - Chicago Lawn
$ws.Range("B66").Value = 2

# Row 80 - Kenwood
$ws.Range("B80").Value = 2
